$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to Text format before writing so Excel
# does not auto-convert numeric-looking strings (e.g. "314.28", "7.340",
# "28.252.92") into floating point numbers, which would lose precision
# and trailing zeros.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '28.252.92'
$ws.Range('E2').Value = '  -0.53%  '
$ws.Range('D3').Value = '1.804.64'
$ws.Range('E3').Value = '  -1.02%  '
$ws.Range('D5').Value = '314.28'
$ws.Range('E5').Value = '  -0.21%  '
$ws.Range('E6').Value = '  +0.44%  '
$ws.Range('D7').Value = '0.5253'
$ws.Range('E7').Value = '  +2.42%  '
$ws.Range('D8').Value = '0.3817'
$ws.Range('E8').Value = '  -2.92%  '
$ws.Range('D9').Value = '0.07913'
$ws.Range('E9').Value = '  +3.26%  '
$ws.Range('D10').Value = '41.77'
$ws.Range('E10').Value = '  +0.37%  '
$ws.Range('D11').Value = '1.097'
$ws.Range('E11').Value = '  -1.10%  '
$ws.Range('D12').Value = '6.322'
$ws.Range('E12').Value = '  +0.85%  '
$ws.Range('E13').Value = '  +0.47%  '
$ws.Range('D14').Value = '20.65'
$ws.Range('E14').Value = '  -1.65%  '
$ws.Range('D15').Value = '1.806.90'
$ws.Range('E15').Value = '  -0.89%  '
$ws.Range('D16').Value = '7.340'
$ws.Range('E16').Value = '  -2.00%  '
$ws.Range('D17').Value = '92.49'
$ws.Range('E17').Value = '  -0.86%  '
$ws.Range('D18').Value = '0.00001090'
$ws.Range('E18').Value = '  -0.57%  '
$ws.Range('E19').Value = '  -0.90%  '
$ws.Range('E20').Value = '  +0.41%  '
$ws.Range('D21').Value = '17.41'
$ws.Range('E21').Value = '  -1.68%  '
$ws.Range('D22').Value = '5.982'
$ws.Range('E22').Value = '  -2.19%  '
$ws.Range('D23').Value = '28.303.65'
$ws.Range('E23').Value = '  -0.39%  '
$ws.Range('D24').Value = '11.15'
$ws.Range('E24').Value = '  -0.25%  '
$ws.Range('D25').Value = '2.237'
$ws.Range('E25').Value = '  -0.85%  '
$ws.Range('D26').Value = '157.52'
$ws.Range('E26').Value = '  +0.66%  '
$ws.Range('D27').Value = '20.52'
$ws.Range('E27').Value = '  -1.26%  '
$ws.Range('D28').Value = '2.012.85'
$ws.Range('E28').Value = '  -1.07%  '
$ws.Range('D29').Value = '2.406'
$ws.Range('E29').Value = '  +0.52%  '
$ws.Range('D30').Value = '123.14'
$ws.Range('E30').Value = '  -0.76%  '
$ws.Range('D31').Value = '0.1104'
$ws.Range('E31').Value = '  +0.93%  '
$ws.Range('D32').Value = '1.060'
$ws.Range('E32').Value = '  -4.39%  '
$ws.Range('E33').Value = '  +0.38%  '
$ws.Range('D34').Value = '5.576'
$ws.Range('E34').Value = '  -1.35%  '
$ws.Range('D35').Value = '0.07193'
$ws.Range('E35').Value = '  +1.38%  '
$ws.Range('D36').Value = '12.11'
$ws.Range('E36').Value = '  +8.03%  '
$ws.Range('D37').Value = '0.2169'
$ws.Range('E37').Value = '  -1.87%  '
$ws.Range('D38').Value = '0.02311'
$ws.Range('E38').Value = '  -0.64%  '
$ws.Range('D39').Value = '8.733'
$ws.Range('E39').Value = '  -0.59%  '
$ws.Range('D40').Value = '5.036'
$ws.Range('E40').Value = '  -2.59%  '
$ws.Range('E41').Value = '  -1.05%  '
$ws.Range('D42').Value = '1.170'
$ws.Range('E42').Value = '  -0.09%  '
$ws.Range('E43').Value = '  -0.66%  '
$ws.Range('D44').Value = '0.6039'
$ws.Range('E44').Value = '  +2.58%  '
$ws.Range('D45').Value = '13.18'
$ws.Range('E45').Value = '  -1.60%  '
$ws.Range('D46').Value = '3.764'
$ws.Range('E46').Value = '  +1.09%  '
$ws.Range('E47').Value = '  +0.82%  '
$ws.Range('E48').Value = '  +0.87%  '
$ws.Range('D49').Value = '1.930'
$ws.Range('E49').Value = '  -2.46%  '
$ws.Range('D50').Value = '0.06826'
$ws.Range('E50').Value = '  -1.07%  '
$ws.Range('D51').Value = '72.91'
$ws.Range('E51').Value = '  -1.52%  '

# Restore the default (unstyled) cell style now that the text values are set,
# matching the original workbook formatting (these cells carry no explicit style).
$dataRange.Style = "Normal"
